# This script applies the weekly update described in the commit:
#   "Fruta / hortaliza, semanal"
#
# It inserts 3 new data rows (515-517) into the "Naranja" price sheet,
# pushing the existing rows 515-548 down to 518-551 (Excel's native
# Insert() shift-down behaviour takes care of that automatically,
# including updating the used-range dimension).
#
# The three new rows carry a new weekly price report for "Navel Late"
# oranges (two quality grades plus an extra packaging unit), dated
# 44826 (2022-09-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at the top of the block (rows 515-517), which
# shifts all rows from 515 downward by 3 (515->518 ... 548->551).
$ws.Rows("515:517").Insert()

function Set-NarRow {
    param($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $origen, $precioKg, $kgUnidad)

    $ws.Cells.Item($row, 1).Value = 4
    $ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($row, 3).Value = "Los Lagos"
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = 10
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102005
    $ws.Cells.Item($row, 10).Value = "Naranja"
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

# New row 515
Set-NarRow 515 44826 "Navel Late" "Primera" 400 15000 16000 15500 "$/caja 15 kilos empedrada" "Región de O'Higgins" 1033 15

# New row 516
Set-NarRow 516 44826 "Navel Late" "Segunda" 200 13000 13000 13000 "$/caja 15 kilos empedrada" "Región de O'Higgins" 867 15

# New row 517
Set-NarRow 517 44826 "Navel Late" "Segunda" 400 8000 8500 8250 "$/malla 16 kilos" "Región de O'Higgins" 516 16
